# The edit rotates the "observation" data (everything that identifies a
# single species record: id, taxon sort order, red-list status, taxon id,
# species name, scientific name, author, coordinates, start/end time and
# the optional public comment) among rows 2-10, while location/date/
# observer metadata (columns that are identical on every row) stay put.
#
# Concretely, new row R receives the old row $mapping[R]'s values for the
# columns below; this reproduces the permutation described by the diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns whose values move between rows (A,B,D,E,F,G,H,Q,R,Z,AB,AC)
$cols = @(1, 2, 4, 5, 6, 7, 8, 17, 18, 26, 28, 29)

# new row -> old row that supplies its data
$mapping = @{
    2  = 4
    3  = 7
    4  = 8
    5  = 2
    6  = 5
    7  = 9
    8  = 10
    9  = 6
    10 = 3
}

# Snapshot the current (pre-edit) values for every relevant cell first,
# so that writes to one row never clobber data still needed for another.
$orig = @{}
foreach ($r in $mapping.Keys) {
    $rowVals = @{}
    foreach ($c in $cols) {
        $rowVals[$c] = $ws.Cells.Item($r, $c).Value()
    }
    $orig[$r] = $rowVals
}

foreach ($r in $mapping.Keys) {
    $src = $mapping[$r]
    foreach ($c in $cols) {
        $val = $orig[$src][$c]
        if ($null -eq $val) {
            $ws.Cells.Item($r, $c).Value = ""
        } else {
            $ws.Cells.Item($r, $c).Value = $val
        }
    }
}
